$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the vector/matrix values in rows 2-15 (columns B:F) to the
# corrected values from the ysrs2 fix.
$ws.Range("B2").Value = 0.2776386982356317
$ws.Range("C2").Value = -0.23604712188397
$ws.Range("D2").Value = -0.4182492285314239
$ws.Range("E2").Value = 0.4142752452009365
$ws.Range("F2").Value = -0.2737950593155914
$ws.Range("B3").Value = 0.2751397658095718
$ws.Range("C3").Value = -0.3223444026074188
$ws.Range("D3").Value = -0.292336359186932
$ws.Range("E3").Value = 0.2923832480224156
$ws.Range("F3").Value = -0.09365034383288866
$ws.Range("B4").Value = 0.2727986977661587
$ws.Range("C4").Value = -0.3611872270012924
$ws.Range("D4").Value = -0.1115849884548746
$ws.Range("E4").Value = 0.04301590552305697
$ws.Range("F4").Value = 0.2898385769221218
$ws.Range("B5").Value = 0.2710486637276452
$ws.Range("C5").Value = -0.3423279971247559
$ws.Range("D5").Value = 0.07194096548077974
$ws.Range("E5").Value = -0.2125818252222675
$ws.Range("F5").Value = 0.4242803674394535
$ws.Range("B6").Value = 0.27050033364696
$ws.Range("C6").Value = -0.2643197397684402
$ws.Range("D6").Value = 0.2165727820829542
$ws.Range("E6").Value = -0.3577389886532019
$ws.Range("F6").Value = 0.1507415991207093
$ws.Range("B7").Value = 0.2698146842081763
$ws.Range("C7").Value = -0.153783131185667
$ws.Range("D7").Value = 0.3208102617733549
$ws.Range("E7").Value = -0.2679692699967574
$ws.Range("F7").Value = -0.2569543648319169
$ws.Range("B8").Value = 0.269017449167217
$ws.Range("C8").Value = -0.02474894780885691
$ws.Range("D8").Value = 0.3572792042086653
$ws.Range("E8").Value = -0.008397440817355857
$ws.Range("F8").Value = -0.4772293225429007
$ws.Range("B9").Value = 0.2681269930992963
$ws.Range("C9").Value = 0.1072068281710587
$ws.Range("D9").Value = 0.3328453965282794
$ws.Range("E9").Value = 0.2259209801083708
$ws.Range("F9").Value = -0.2736441462810427
$ws.Range("B10").Value = 0.2667137845380258
$ws.Range("C10").Value = 0.2224814403515821
$ws.Range("D10").Value = 0.2713168401051604
$ws.Range("E10").Value = 0.3523462173661724
$ws.Range("F10").Value = 0.185603807348223
$ws.Range("B11").Value = 0.2652940227031764
$ws.Range("C11").Value = 0.3115063811691925
$ws.Range("D11").Value = 0.1303826386438549
$ws.Range("E11").Value = 0.2566529246012451
$ws.Range("F11").Value = 0.3899834922421908
$ws.Range("B12").Value = 0.2637765907961352
$ws.Range("C12").Value = 0.3541382618721791
$ws.Range("D12").Value = -0.05961637615235373
$ws.Range("E12").Value = 0.04156936701049795
$ws.Range("F12").Value = 0.2079904842033242
$ws.Range("B13").Value = 0.2610636829777343
$ws.Range("C13").Value = 0.3355258761244662
$ws.Range("D13").Value = -0.2138283405220615
$ws.Range("E13").Value = -0.1564711046603464
$ws.Range("F13").Value = -0.02722959222535455
$ws.Range("B14").Value = 0.2568204368805638
$ws.Range("C14").Value = 0.2655545728250846
$ws.Range("D14").Value = -0.2965068712464946
$ws.Range("E14").Value = -0.2966815092273329
$ws.Range("F14").Value = -0.1185747259579901
$ws.Range("B15").Value = 0.2527632363148538
$ws.Range("C15").Value = 0.179070408140679
$ws.Range("D15").Value = -0.3254495131697182
$ws.Range("E15").Value = -0.3751018790046719
$ws.Range("F15").Value = -0.1275569361336337

# Remove the now-obsolete trailing rows 16 and 17 (indices 14 and 15),
# shrinking the sheet dimension from A1:F17 to A1:F15.
$ws.Range("A16:F17").EntireRow.Delete()
